$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing test-case row's name to the "_1" variant (shared
# string "Apache_POI_TC" -> "Apache_POI_TC_1") - this is the data used by
# the new DataProvider-driven test class.
$ws.Range("A2").Value = "Apache_POI_TC_1"

# Add a brand new data row for the second DataProvider entry.
$ws.Range("A3").Value = "Apache_POI_TC_2"
$ws.Range("B3").Value = "testuser_2"
$ws.Range("C3").Value = "Test@1234"

# New hyperlink for the new password cell, mirroring the one on C2.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Test@1234")
$ws.Range("C3").Style = "Hyperlink"

# Match the saved selection from the edited workbook.
[void]$ws.Range("D3").Select()
